# Updated cryptos list (price / 1h-volume refresh), mirrors the GitHub
# Actions bot commit. Most cells are plain text re-assignments; a handful
# of Price cells (column D) whose new text happens to look like a bare
# number ("215.90", "9.50", ...) are round-tripped through a temporary
# "@" (Text) number format so Excel stores them as strings instead of
# silently coercing them to doubles (which would also eat the trailing
# zero). The original style is restored immediately after so no cell
# ends up with a different style than it started with.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.723.93'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '1.647.21'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  +0.13%  '
$origStyle_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.90'
$ws.Range('D5').Style = $origStyle_D5
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('E6').Value = '  +1.59%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  +1.56%  '
$origStyle_D9 = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0627'
$ws.Range('D9').Style = $origStyle_D9
$ws.Range('E9').Value = '  +0.65%  '
$origStyle_D10 = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.24'
$ws.Range('D10').Style = $origStyle_D10
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').Value = '1.876.51'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.632.32'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$origStyle_D14 = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.18'
$ws.Range('D14').Style = $origStyle_D14
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('E15').Value = '  +2.03%  '
$origStyle_D16 = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.43'
$ws.Range('D16').Style = $origStyle_D16
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = '26.722.63'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('E18').Value = '  +0.56%  '
$origStyle_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '218.80'
$ws.Range('D19').Style = $origStyle_D19
$ws.Range('E19').Value = '  +2.19%  '
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('E21').Value = '  +1.56%  '
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +14.13%  '
$origStyle_D24 = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.50'
$ws.Range('D24').Style = $origStyle_D24
$ws.Range('E24').Value = '  +2.38%  '
$origStyle_D25 = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.97'
$ws.Range('D25').Style = $origStyle_D25
$ws.Range('E25').Value = '  -1.80%  '
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E28').Value = '  +4.41%  '
$origStyle_D29 = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.73'
$ws.Range('D29').Style = $origStyle_D29
$ws.Range('E29').Value = '  +1.40%  '
$origStyle_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0517'
$ws.Range('D30').Style = $origStyle_D30
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('D34').Value = '1.282.91'
$ws.Range('E34').Value = '  +5.43%  '
$ws.Range('E35').Value = '  +3.86%  '
$origStyle_D36 = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.43'
$ws.Range('D36').Style = $origStyle_D36
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('E37').Value = '  +3.10%  '
$origStyle_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.536'
$ws.Range('D38').Style = $origStyle_D38
$ws.Range('E38').Value = '  +6.15%  '
$ws.Range('E39').Value = '  +4.05%  '
$ws.Range('E40').Value = '  +0.19%  '
$origStyle_D41 = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.814'
$ws.Range('D41').Style = $origStyle_D41
$ws.Range('E41').Value = '  +2.86%  '
$origStyle_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.26'
$ws.Range('D42').Style = $origStyle_D42
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('E43').Value = '  +1.86%  '
$ws.Range('D44').Value = '1.786.73'
$ws.Range('E44').Value = '  +1.28%  '
$origStyle_D45 = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '92.08'
$ws.Range('D45').Style = $origStyle_D45
$ws.Range('E45').Value = '  -0.78%  '
$origStyle_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '59.92'
$ws.Range('D46').Style = $origStyle_D46
$ws.Range('E46').Value = '  +9.32%  '
$ws.Range('E47').Value = '  +1.60%  '
$ws.Range('E48').Value = '  +1.18%  '
$origStyle_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.75'
$ws.Range('D49').Style = $origStyle_D49
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('E50').Value = '  +1.94%  '
$ws.Range('E51').Value = '  -0.03%  '
